$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 217/218 header-like text cells (order matters for shared-string index allocation) ---
$ws.Range("K217").Value = "PQI"
$ws.Range("K218").Value = "Unit"
$ws.Range("L218").Value = "Inst"
$ws.Range("M218").Value = "Dev"
$ws.Range("N218").Value = "#Automated/#dev"
$ws.Range("O218").Value = "Output comp"
$ws.Range("P218").Value = "Total"
$ws.Range("N217").Value = "more better"
$ws.Range("O217").Value = "more better"
$ws.Range("P217").Value = "more better"

# --- K/L/M/O numeric data for rows 219-248 ---
$ws.Range("K219").Value = 129
$ws.Range("L219").Value = 5
$ws.Range("M219").Value = 8
$ws.Range("O219").Value = 1
$ws.Range("K220").Value = 53
$ws.Range("L220").Value = 19
$ws.Range("M220").Value = 4
$ws.Range("O220").Value = 1
$ws.Range("K221").Value = 21
$ws.Range("L221").Value = 2
$ws.Range("M221").Value = 2
$ws.Range("O221").Value = 0.5
$ws.Range("K222").Value = 134
$ws.Range("L222").Value = 4
$ws.Range("M222").Value = 3
$ws.Range("O222").Value = 1
$ws.Range("K223").Value = 125
$ws.Range("L223").Value = 1
$ws.Range("M223").Value = 3
$ws.Range("O223").Value = 1
$ws.Range("K224").Value = 240
$ws.Range("L224").Value = 8
$ws.Range("M224").Value = 4
$ws.Range("O224").Value = 1
$ws.Range("K225").Value = 300
$ws.Range("L225").Value = 2
$ws.Range("M225").Value = 2
$ws.Range("O225").Value = 1
$ws.Range("K226").Value = 14
$ws.Range("L226").Value = 1
$ws.Range("M226").Value = 2
$ws.Range("O226").Value = 1
$ws.Range("K227").Value = 37
$ws.Range("L227").Value = 2
$ws.Range("M227").Value = 2
$ws.Range("O227").Value = 1
$ws.Range("K228").Value = 40
$ws.Range("L228").Value = 6
$ws.Range("M228").Value = 4
$ws.Range("O228").Value = 1
$ws.Range("K229").Value = 50
$ws.Range("L229").Value = 1
$ws.Range("M229").Value = 3
$ws.Range("O229").Value = 0.5
$ws.Range("K230").Value = 39
$ws.Range("L230").Value = 1
$ws.Range("M230").Value = 2
$ws.Range("O230").Value = 1
$ws.Range("K231").Value = 600
$ws.Range("L231").Value = 25
$ws.Range("M231").Value = 0
$ws.Range("O231").Value = 1
$ws.Range("K232").Value = 35
$ws.Range("L232").Value = 5
$ws.Range("M232").Value = 20
$ws.Range("O232").Value = 1
$ws.Range("K233").Value = 61
$ws.Range("L233").Value = 2
$ws.Range("M233").Value = 3
$ws.Range("O233").Value = 1
$ws.Range("K234").Value = 23
$ws.Range("L234").Value = 1
$ws.Range("M234").Value = 4
$ws.Range("O234").Value = 1
$ws.Range("K235").Value = 14
$ws.Range("L235").Value = 1
$ws.Range("M235").Value = 3
$ws.Range("O235").Value = 1
$ws.Range("K236").Value = 70
$ws.Range("L236").Value = 2
$ws.Range("M236").Value = 0
$ws.Range("O236").Value = 1
$ws.Range("K237").Value = 17
$ws.Range("L237").Value = 1
$ws.Range("M237").Value = 1
$ws.Range("O237").Value = 1
$ws.Range("K238").Value = 50
$ws.Range("L238").Value = 2
$ws.Range("M238").Value = 1
$ws.Range("O238").Value = 1
$ws.Range("K239").Value = 90
$ws.Range("L239").Value = 1
$ws.Range("M239").Value = 2
$ws.Range("O239").Value = 1
$ws.Range("K240").Value = 28
$ws.Range("L240").Value = 1
$ws.Range("M240").Value = 5
$ws.Range("O240").Value = 1
$ws.Range("K241").Value = 23
$ws.Range("L241").Value = 1
$ws.Range("M241").Value = 1
$ws.Range("O241").Value = 1
$ws.Range("K242").Value = 58
$ws.Range("L242").Value = 6
$ws.Range("M242").Value = 4
$ws.Range("O242").Value = 1
$ws.Range("K243").Value = 165
$ws.Range("L243").Value = 5
$ws.Range("M243").Value = 1
$ws.Range("O243").Value = 1
$ws.Range("K244").Value = 200
$ws.Range("L244").Value = 4
$ws.Range("M244").Value = 9
$ws.Range("O244").Value = 1
$ws.Range("K245").Value = 80
$ws.Range("L245").Value = 2
$ws.Range("M245").Value = 0
$ws.Range("O245").Value = 1
$ws.Range("K246").Value = 71
$ws.Range("L246").Value = 5
$ws.Range("M246").Value = 2
$ws.Range("O246").Value = 1
$ws.Range("K247").Value = 68
$ws.Range("L247").Value = 1
$ws.Range("M247").Value = 3
$ws.Range("O247").Value = 1
$ws.Range("K248").Value = 33
$ws.Range("L248").Value = 2
$ws.Range("M248").Value = 0
$ws.Range("O248").Value = 1

# --- N column formulas: N219 individual, N220:N248 one shared group ---
$ws.Range("N219").Formula = "=(K219+M219)/(K219+L219)"
$ws.Range("N220:N248").Formula = "=(K220+M220)/(K220+L220)"

# --- P column formulas: P219 & P220 individual, then contiguous shared blocks between
#     the manually re-typed "break" cells (226,228,234,235,241,243,247), which are individual too ---
$ws.Range("P219").Formula = "=N219 * O219"
$ws.Range("P220").Formula = "=N220 * O220"
$ws.Range("P221:P225").Formula = "=N221 * O221"
$ws.Range("P226").Formula = "=N226 * O226"
$ws.Range("P227").Formula = "=N227 * O227"
$ws.Range("P228").Formula = "=N228 * O228"
$ws.Range("P229:P233").Formula = "=N229 * O229"
$ws.Range("P234").Formula = "=N234 * O234"
$ws.Range("P235").Formula = "=N235 * O235"
$ws.Range("P236:P240").Formula = "=N236 * O236"
$ws.Range("P241").Formula = "=N241 * O241"
$ws.Range("P242").Formula = "=N242 * O242"
$ws.Range("P243").Formula = "=N243 * O243"
$ws.Range("P244:P246").Formula = "=N244 * O244"
$ws.Range("P247").Formula = "=N247 * O247"
$ws.Range("P248").Formula = "=N248 * O248"

# --- Column widths for N, O, P (best achievable given engine quantization) ---
$ws.Columns.Item(14).ColumnWidth = 16.85
$ws.Columns.Item(15).ColumnWidth = 11.66
$ws.Columns.Item(16).ColumnWidth = 13.17

# --- Sheet view: zoom, selection, scroll position ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("B228").Select()
$excel.ActiveWindow.ScrollRow = 228
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("P248").Select()
